$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = '{''22'', ''240'', ''245'', ''1'', ''37'', ''41'', ''40'', ''3''}'
$ws.Range("N3").Value = '{''245'', ''37'', ''1'', ''41'', ''3''}'
$ws.Range("N4").Value = '{''245'', ''37'', ''1'', ''41'', ''3''}'
$ws.Range("N6").Value = '{''19'', ''7'', ''1''}'
$ws.Range("N7").Value = '{''54'', ''53'', ''19'', ''7'', ''1'', ''246'', ''51''}'
$ws.Range("N8").Value = '{''19'', ''4'', ''7'', ''1''}'
$ws.Range("N9").Value = '{''24'', ''19'', ''7'', ''1'', ''52''}'
$ws.Range("N10").Value = '{''19'', ''7'', ''1''}'
$ws.Range("N11").Value = '{''19'', ''7'', ''1'', ''5'', ''12'', ''2''}'
$ws.Range("N15").Value = '{''247'', ''15'', ''17'', ''1''}'
$ws.Range("N17").Value = '{''50'', ''51'', ''17'', ''1''}'
$ws.Range("N18").Value = '{''63'', ''17'', ''1'', ''18'', ''2''}'
$ws.Range("N35").Value = '{''63'', ''17'', ''1'', ''18'', ''2''}'
$ws.Range("N36").Value = '{''50'', ''51'', ''17'', ''1''}'
$ws.Range("N38").Value = '{''247'', ''15'', ''17'', ''1''}'
$ws.Range("N42").Value = '{''19'', ''7'', ''1'', ''5'', ''12'', ''2''}'
$ws.Range("N43").Value = '{''19'', ''7'', ''1''}'
$ws.Range("N44").Value = '{''24'', ''19'', ''7'', ''1'', ''52''}'
$ws.Range("N45").Value = '{''19'', ''4'', ''7'', ''1''}'
$ws.Range("N46").Value = '{''54'', ''53'', ''19'', ''7'', ''1'', ''246'', ''51''}'
$ws.Range("N47").Value = '{''19'', ''7'', ''1''}'
$ws.Range("N49").Value = '{''245'', ''37'', ''1'', ''41'', ''3''}'
$ws.Range("N50").Value = '{''245'', ''37'', ''1'', ''41'', ''3''}'
$ws.Range("N51").Value = '{''22'', ''240'', ''245'', ''1'', ''37'', ''41'', ''40'', ''3''}'
$ws.Range("N52").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N53").Value = '{''12'', ''2'', ''17'', ''13''}'
$ws.Range("N54").Value = '{''14'', ''24'', ''13'', ''17'', ''4'', ''12'', ''2''}'
$ws.Range("N55").Value = '{''12'', ''2''}'
$ws.Range("N56").Value = '{''12'', ''2''}'
$ws.Range("N57").Value = '{''12'', ''4'', ''2''}'
$ws.Range("N58").Value = '{''12'', ''4'', ''2''}'
$ws.Range("N59").Value = '{''12'', ''5'', ''2''}'
$ws.Range("N60").Value = '{''12'', ''5'', ''3''}'
$ws.Range("N61").Value = '{''12'', ''5'', ''3'', ''24''}'
$ws.Range("N62").Value = '{''12'', ''3'', ''24''}'
$ws.Range("N63").Value = '{''12'', ''246''}'
$ws.Range("N64").Value = '{''12'', ''246'', ''65''}'
$ws.Range("N65").Value = '{''245'', ''62'', ''4'', ''65'', ''12''}'
$ws.Range("N66").Value = '{''245'', ''62'', ''4'', ''65'', ''12''}'
$ws.Range("N67").Value = '{''245'', ''62'', ''4'', ''65'', ''12''}'
$ws.Range("N68").Value = '{''12'', ''245'', ''65'', ''62''}'
$ws.Range("N69").Value = '{''12'', ''54'', ''51'', ''53''}'
$ws.Range("N70").Value = '{''12'', ''245'', ''65'', ''62''}'
$ws.Range("N71").Value = '{''245'', ''62'', ''4'', ''65'', ''12''}'
$ws.Range("N72").Value = '{''245'', ''62'', ''4'', ''65'', ''12''}'
$ws.Range("N73").Value = '{''245'', ''62'', ''4'', ''65'', ''12''}'
$ws.Range("N74").Value = '{''12'', ''246'', ''65''}'
$ws.Range("N75").Value = '{''12'', ''246''}'
$ws.Range("N76").Value = '{''12'', ''3'', ''24''}'
$ws.Range("N77").Value = '{''12'', ''5'', ''3'', ''24''}'
$ws.Range("N78").Value = '{''12'', ''5'', ''3''}'
$ws.Range("N79").Value = '{''12'', ''5'', ''2''}'
$ws.Range("N80").Value = '{''19'', ''7'', ''1'', ''5'', ''12'', ''2''}'
$ws.Range("N81").Value = '{''12'', ''4'', ''2''}'
$ws.Range("N82").Value = '{''12'', ''4'', ''2''}'
$ws.Range("N83").Value = '{''12'', ''2''}'
$ws.Range("N84").Value = '{''12'', ''2''}'
$ws.Range("N85").Value = '{''12'', ''2'', ''17'', ''13''}'
$ws.Range("N86").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N87").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N90").Value = '{''19'', ''7'', ''17'', ''13''}'
$ws.Range("N92").Value = '{''19'', ''18'', ''247'', ''13''}'
$ws.Range("N93").Value = '{''247'', ''18'', ''13''}'
$ws.Range("N94").Value = '{''13'', ''247'', ''15'', ''7'', ''18''}'
$ws.Range("N98").Value = '{''13'', ''19'', ''247'', ''15'', ''7'', ''18'', ''51'', ''50''}'
$ws.Range("N109").Value = '{''13'', ''19'', ''247'', ''15'', ''7'', ''18'', ''51'', ''50''}'
$ws.Range("N113").Value = '{''13'', ''247'', ''15'', ''7'', ''18''}'
$ws.Range("N114").Value = '{''247'', ''18'', ''13''}'
$ws.Range("N115").Value = '{''19'', ''18'', ''247'', ''13''}'
$ws.Range("N117").Value = '{''19'', ''7'', ''17'', ''13''}'
$ws.Range("N120").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N121").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N122").Value = '{''14'', ''24'', ''13'', ''17'', ''4'', ''12'', ''2''}'
$ws.Range("N123").Value = '{''14'', ''52'', ''4'', ''24''}'
$ws.Range("N125").Value = '{''14'', ''54'', ''51'', ''53''}'
$ws.Range("N129").Value = '{''19'', ''14'', ''7''}'
$ws.Range("N132").Value = '{''22'', ''14'', ''65''}'
$ws.Range("N133").Value = '{''22'', ''14'', ''240'', ''245'', ''37'', ''40'', ''3''}'
$ws.Range("N134").Value = '{''14'', ''3'', ''245'', ''37''}'
$ws.Range("N135").Value = '{''14'', ''3'', ''245'', ''37''}'
$ws.Range("N138").Value = '{''14'', ''3'', ''245'', ''37''}'
$ws.Range("N139").Value = '{''14'', ''3'', ''245'', ''37''}'
$ws.Range("N140").Value = '{''22'', ''14'', ''240'', ''245'', ''37'', ''40'', ''3''}'
$ws.Range("N141").Value = '{''22'', ''14'', ''65''}'
$ws.Range("N144").Value = '{''19'', ''14'', ''7''}'
$ws.Range("N148").Value = '{''14'', ''54'', ''51'', ''53''}'
$ws.Range("N150").Value = '{''14'', ''52'', ''4'', ''24''}'
$ws.Range("N151").Value = '{''14'', ''24'', ''13'', ''17'', ''4'', ''12'', ''2''}'
$ws.Range("N152").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N153").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N154").Value = '{''19'', ''5'', ''7'', ''17''}'
$ws.Range("N155").Value = '{''19'', ''7'', ''17'', ''13''}'
$ws.Range("N158").Value = '{''247'', ''15'', ''7'', ''17'', ''18''}'
$ws.Range("N159").Value = '{''247'', ''15'', ''17''}'
$ws.Range("N160").Value = '{''369'', ''63'', ''17''}'
$ws.Range("N169").Value = '{''369'', ''63'', ''17''}'
$ws.Range("N170").Value = '{''247'', ''15'', ''17''}'
$ws.Range("N171").Value = '{''247'', ''15'', ''7'', ''17'', ''18''}'
$ws.Range("N174").Value = '{''19'', ''7'', ''17'', ''13''}'
$ws.Range("N175").Value = '{''5'', ''17''}'
$ws.Range("N176").Value = '{''19'', ''5'', ''7'', ''17''}'
$ws.Range("N177").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N178").Value = '{''22'', ''231'', ''19'', ''369'', ''15'', ''61'', ''51'', ''50'', ''36''}'
$ws.Range("N181").Value = '{''19'', ''247'', ''21''}'
$ws.Range("N183").Value = '{''13'', ''19'', ''247'', ''15'', ''7'', ''18'', ''51'', ''50''}'
$ws.Range("N185").Value = '{''19'', ''5'', ''7'', ''17''}'
$ws.Range("N187").Value = '{''19'', ''41'', ''245'', ''37''}'
$ws.Range("N188").Value = '{''19'', ''41''}'
$ws.Range("N189").Value = '{''19'', ''41''}'
$ws.Range("N190").Value = '{''19'', ''41''}'
$ws.Range("N191").Value = '{''19'', ''40'', ''240'', ''41''}'
$ws.Range("N199").Value = '{''19'', ''40'', ''240'', ''41''}'
$ws.Range("N200").Value = '{''19'', ''41''}'
$ws.Range("N201").Value = '{''19'', ''41''}'
$ws.Range("N202").Value = '{''19'', ''41''}'
$ws.Range("N203").Value = '{''19'', ''41'', ''245'', ''37''}'
$ws.Range("N205").Value = '{''19'', ''7'', ''1'', ''5'', ''12'', ''2''}'
$ws.Range("N207").Value = '{''13'', ''19'', ''247'', ''15'', ''7'', ''18'', ''51'', ''50''}'
$ws.Range("N209").Value = '{''19'', ''247'', ''21''}'
$ws.Range("N212").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N215").Value = '{''24'', ''15'', ''62'', ''51'', ''50'', ''2''}'
$ws.Range("N216").Value = '{''15'', ''2'', ''62''}'
$ws.Range("N219").Value = '{''50'', ''51'', ''2''}'
$ws.Range("N220").Value = '{''63'', ''17'', ''1'', ''18'', ''2''}'
$ws.Range("N222").Value = '{''247'', ''369'', ''2''}'
$ws.Range("N227").Value = '{''247'', ''369'', ''2''}'
$ws.Range("N229").Value = '{''63'', ''17'', ''1'', ''18'', ''2''}'
$ws.Range("N230").Value = '{''50'', ''51'', ''2''}'
$ws.Range("N233").Value = '{''15'', ''2'', ''62''}'
$ws.Range("N234").Value = '{''24'', ''15'', ''62'', ''51'', ''50'', ''2''}'
$ws.Range("N237").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N239").Value = '{''24'', ''19'', ''7'', ''1'', ''52''}'
$ws.Range("N242").Value = '{''246'', ''5'', ''24'', ''65''}'
$ws.Range("N248").Value = '{''24'', ''15'', ''62'', ''51'', ''50'', ''2''}'
$ws.Range("N250").Value = '{''246'', ''461'', ''24'', ''62''}'
$ws.Range("N251").Value = '{''246'', ''461'', ''24'', ''62''}'
$ws.Range("N253").Value = '{''24'', ''15'', ''62'', ''51'', ''50'', ''2''}'
$ws.Range("N259").Value = '{''246'', ''24'', ''65''}'
$ws.Range("N261").Value = '{''12'', ''5'', ''3'', ''24''}'
$ws.Range("N263").Value = '{''24'', ''19'', ''7'', ''1'', ''52''}'
$ws.Range("N265").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N278").Value = '{''463'', ''25'', ''5''}'
$ws.Range("N279").Value = '{''463'', ''25'', ''5'', ''461''}'
$ws.Range("N281").Value = '{''25'', ''15'', ''62'', ''5'', ''65''}'
$ws.Range("N282").Value = '{''25'', ''15'', ''62'', ''5'', ''65''}'
$ws.Range("N284").Value = '{''463'', ''25'', ''5'', ''461''}'
$ws.Range("N285").Value = '{''463'', ''25'', ''5''}'
$ws.Range("N304").Value = '{''26'', ''245'', ''37''}'
$ws.Range("N305").Value = '{''26'', ''7''}'
$ws.Range("N306").Value = '{''26'', ''246'', ''43''}'
$ws.Range("N308").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N310").Value = '{''26'', ''246'', ''43''}'
$ws.Range("N311").Value = '{''26'', ''7''}'
$ws.Range("N312").Value = '{''26'', ''245'', ''37''}'
$ws.Range("N331").Value = '{''22'', ''240'', ''245'', ''1'', ''37'', ''41'', ''40'', ''3''}'
$ws.Range("N346").Value = '{''52'', ''4'', ''62''}'
$ws.Range("N347").Value = '{''52'', ''4'', ''245'', ''65''}'
$ws.Range("N348").Value = '{''245'', ''62'', ''4'', ''65'', ''12''}'
$ws.Range("N349").Value = '{''245'', ''62'', ''4'', ''65'', ''12''}'
$ws.Range("N350").Value = '{''245'', ''62'', ''4'', ''65'', ''12''}'
$ws.Range("N353").Value = '{''246'', ''4''}'
$ws.Range("N355").Value = '{''19'', ''4'', ''7'', ''1''}'
$ws.Range("N356").Value = '{''12'', ''4'', ''2''}'
$ws.Range("N357").Value = '{''12'', ''4'', ''2''}'
$ws.Range("N358").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N359").Value = '{''14'', ''54'', ''24'', ''53'', ''13'', ''17'', ''4'', ''51'', ''12'', ''26'', ''52'', ''2''}'
$ws.Range("N360").Value = '{''12'', ''4'', ''2''}'
$ws.Range("N361").Value = '{''12'', ''4'', ''2''}'
$ws.Range("N362").Value = '{''19'', ''4'', ''7'', ''1''}'
$ws.Range("N363").Value = '{''246'', ''4''}'
$ws.Range("N366").Value = '{''245'', ''62'', ''4'', ''65'', ''12''}'
$ws.Range("N367").Value = '{''245'', ''62'', ''4'', ''65'', ''12''}'
$ws.Range("N368").Value = '{''245'', ''62'', ''4'', ''65'', ''12''}'
$ws.Range("N369").Value = '{''52'', ''4'', ''245'', ''65''}'
$ws.Range("N370").Value = '{''52'', ''4'', ''62''}'
$ws.Range("N373").Value = '{''246'', ''15'', ''5'', ''65''}'
$ws.Range("N374").Value = '{''246'', ''5'', ''24'', ''65''}'
$ws.Range("N382").Value = '{''5'', ''17''}'
$ws.Range("N383").Value = '{''246'', ''15'', ''5'', ''65''}'
$ws.Range("N385").Value = '{''246'', ''43'', ''7'', ''65''}'
$ws.Range("N386").Value = '{''246'', ''43'', ''7'', ''65''}'
$ws.Range("N387").Value = '{''26'', ''7''}'
$ws.Range("N389").Value = '{''22'', ''7''}'
$ws.Range("N390").Value = '{''22'', ''7''}'
$ws.Range("N391").Value = '{''19'', ''14'', ''7''}'
$ws.Range("N392").Value = '{''247'', ''15'', ''7'', ''17'', ''18''}'
$ws.Range("N393").Value = '{''247'', ''18'', ''15'', ''7''}'
$ws.Range("N394").Value = '{''13'', ''247'', ''15'', ''7'', ''18''}'
$ws.Range("N395").Value = '{''13'', ''19'', ''247'', ''15'', ''7'', ''18'', ''51'', ''50''}'
$ws.Range("N396").Value = '{''15'', ''7''}'
$ws.Range("N397").Value = '{''247'', ''15'', ''21'', ''7''}'
$ws.Range("N398").Value = '{''50'', ''51'', ''21'', ''7''}'
$ws.Range("N399").Value = '{''21'', ''7''}'
$ws.Range("N400").Value = '{''21'', ''7''}'
$ws.Range("N401").Value = '{''369'', ''21'', ''7''}'
$ws.Range("N402").Value = '{''369'', ''7''}'
$ws.Range("N404").Value = '{''369'', ''7''}'
$ws.Range("N405").Value = '{''369'', ''21'', ''7''}'
$ws.Range("N406").Value = '{''21'', ''7''}'
$ws.Range("N407").Value = '{''21'', ''7''}'
$ws.Range("N408").Value = '{''50'', ''51'', ''21'', ''7''}'
$ws.Range("N409").Value = '{''247'', ''15'', ''21'', ''7''}'
$ws.Range("N410").Value = '{''13'', ''19'', ''247'', ''15'', ''7'', ''18'', ''51'', ''50''}'
$ws.Range("N411").Value = '{''13'', ''247'', ''15'', ''7'', ''18''}'
$ws.Range("N412").Value = '{''247'', ''18'', ''15'', ''7''}'
$ws.Range("N413").Value = '{''19'', ''14'', ''7''}'
$ws.Range("N414").Value = '{''22'', ''7''}'
$ws.Range("N415").Value = '{''22'', ''7''}'
$ws.Range("N417").Value = '{''26'', ''7''}'
$ws.Range("N418").Value = '{''246'', ''43'', ''7'', ''65''}'
$ws.Range("N419").Value = '{''246'', ''43'', ''7'', ''65''}'
